# Apply cryptos list update (Sat Nov 18 02:41:06 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.468.17"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").Value = "1.953.40"
$ws.Range("E3").Value = "  -1.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'244.00"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6
$ws.Range("E6").Value = "  -1.92%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.60"
$ws.Range("E8").Value = "  -2.30%  "

# Row 9
$ws.Range("D9").Value = "'0.365"
$ws.Range("E9").Value = "  -3.23%  "

# Row 10
$ws.Range("D10").Value = "'0.0852"
$ws.Range("E10").Value = "  +3.86%  "

# Row 11
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
$ws.Range("D12").Value = "2.241.19"
$ws.Range("E12").Value = "  -1.74%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.821"
$ws.Range("E13").Value = "  -5.64%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.46"
$ws.Range("E14").Value = "  -12.98%  "

# Row 15
$ws.Range("D15").Value = "'13.54"
$ws.Range("E15").Value = "  -4.59%  "

# Row 16
$ws.Range("E16").Value = "  -4.64%  "

# Row 17
$ws.Range("D17").Value = "1.956.04"
$ws.Range("E17").Value = "  -1.74%  "

# Row 18
$ws.Range("D18").Value = "36.378.09"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0884"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("D20").Value = "'69.80"
$ws.Range("E20").Value = "  -2.18%  "

# Row 21
$ws.Range("D21").Value = "'229.96"
$ws.Range("E21").Value = "  -2.40%  "

# Row 22
$ws.Range("D22").Value = "'5.07"
$ws.Range("E22").Value = "  -5.00%  "

# Row 23
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").Value = "'2.42"
$ws.Range("E24").Value = "  -8.62%  "

# Row 25
$ws.Range("E25").Value = "  -0.73%  "

# Row 26
$ws.Range("D26").Value = "'9.27"
$ws.Range("E26").Value = "  -8.24%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'161.56"
$ws.Range("E27").Value = "  -0.38%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  +4.58%  "

# Row 29
$ws.Range("D29").Value = "'19.45"
$ws.Range("E29").Value = "  -2.36%  "

# Row 30
$ws.Range("D30").Value = "'0.118"
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  -1.92%  "

# Row 32
$ws.Range("D32").Value = "'4.66"
$ws.Range("E32").Value = "  -5.71%  "

# Row 33
$ws.Range("D33").Value = "'0.0652"
$ws.Range("E33").Value = "  +3.41%  "

# Row 34
$ws.Range("D34").Value = "'4.28"
$ws.Range("E34").Value = "  -4.54%  "

# Row 35
$ws.Range("D35").Value = "'6.19"
$ws.Range("E35").Value = "  -0.74%  "

# Row 37
$ws.Range("E37").Value = "  +1.24%  "

# Row 38
$ws.Range("E38").Value = "  -6.18%  "

# Row 39
$ws.Range("E39").Value = "  -2.32%  "

# Row 40
$ws.Range("D40").Value = "'0.0985"
$ws.Range("E40").Value = "  +0.97%  "

# Row 41
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("E42").Value = "  -6.44%  "

# Row 43
$ws.Range("E43").Value = "  -1.52%  "

# Row 44
$ws.Range("D44").Value = "'15.75"
$ws.Range("E44").Value = "  -3.87%  "

# Row 45
$ws.Range("D45").Value = "1.358.41"
$ws.Range("E45").Value = "  -1.13%  "

# Row 46
$ws.Range("E46").Value = "  -6.69%  "

# Row 47
$ws.Range("D47").Value = "'87.90"
$ws.Range("E47").Value = "  -5.58%  "

# Row 48
$ws.Range("D48").Value = "'7.16"
$ws.Range("E48").Value = "  -6.69%  "

# Row 49
$ws.Range("E49").Value = "  -0.71%  "

# Row 50
$ws.Range("D50").Value = "'44.99"
$ws.Range("E50").Value = "  -1.35%  "

# Row 51
$ws.Range("D51").Value = "2.131.56"
$ws.Range("E51").Value = "  -2.00%  "

